$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.430.86'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '2.649.37'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.04'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.97%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.143'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.31%  '
$ws.Range("B10").Value = 'TRON'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.157'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.27'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.353'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.15'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.25%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.132.28'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000188'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '68.364.96'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.647.12'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.64'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.69%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '364.48'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.48'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.42'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.78%  '
$ws.Range("B22").Value = 'NEARProtocol'
$ws.Range("C22").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.78'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.09'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.08%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.09'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.97'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.35%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '2.799.94'
$ws.Range("E27").Value = '  +0.90%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000103'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.69%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '572.96'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.37%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.07'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.41'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.79%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.87'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.03%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.64'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.46%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.129'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.67%  '
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.77'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.65'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.99%  '
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.371'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.88'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.33'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.65'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("B43").Value = 'BabyDogeCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D43").Value = '0.0₆0321'
$ws.Range("E43").Value = '  -4.88%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '158.55'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.82'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.31%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.92'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("B48").Value = 'Optimism'
$ws.Range("C48").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.70'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0779'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.73%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.575'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.69%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.616'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.11%  '
